# Germany Landesliga workbook update (19-06-2024 21:51)
#
# The underlying data feed re-sorted two pairs of match rows:
#   - row 12  <->  row 13   (match ids 7035046 / 7035047)
#   - row 129 <->  row 130  (match ids 8271342 / 8271343)
#
# In each pair, every column except the running "id" in column A
# (which stays tied to its row position) swaps between the two rows.
# The HomeTeam / AwayTeam names referenced elsewhere in the sheet keep
# their same text, so no further edits are required there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastCol = 30   # column AD

# --- swap row 12 and row 13 (keep column A as-is) ---
$row1 = 12
$row2 = 13
for ($c = 2; $c -le $lastCol; $c++) {
    $cell1 = $ws.Cells.Item($row1, $c)
    $cell2 = $ws.Cells.Item($row2, $c)
    $v1 = $cell1.Value2
    $v2 = $cell2.Value2
    $cell1.Value = $v2
    $cell2.Value = $v1
}

# --- swap row 129 and row 130 (keep column A as-is) ---
$row1 = 129
$row2 = 130
for ($c = 2; $c -le $lastCol; $c++) {
    $cell1 = $ws.Cells.Item($row1, $c)
    $cell2 = $ws.Cells.Item($row2, $c)
    $v1 = $cell1.Value2
    $v2 = $cell2.Value2
    $cell1.Value = $v2
    $cell2.Value = $v1
}

Write-Output "Rows 12/13 and 129/130 swapped"
